# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet to the refreshed values scraped by the GitHub Actions job.
#
# Price/volume cells in this sheet are stored as text (e.g. "213.03",
# "  -0.35%  ") rather than numbers, since some prices use "." as a
# thousands separator (e.g. "28.539.53"). Several of the new values look
# like ordinary decimal numbers to Excel (e.g. "213.01", "63.00"), so a
# plain Value assignment would silently convert them into numeric cells
# and normalize/round their text (dropping the trailing zero in "63.00",
# etc). To keep them as text - matching the original cell type and exact
# string - the number format is forced to "@" (Text) before assigning the
# value, and the style is reset back to Normal afterwards so no visible
# formatting change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# row 2 - Bitcoin
Set-TextValue "D2" "28.558.77"
$ws.Range("E2").Value = "  +1.05%  "

# row 3 - Ethereum
Set-TextValue "D3" "1.572.48"
$ws.Range("E3").Value = "  -1.15%  "

# row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# row 5 - BNB
Set-TextValue "D5" "213.01"
$ws.Range("E5").Value = "  -0.33%  "

# row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# row 8 - OKB
Set-TextValue "D8" "45.60"
$ws.Range("E8").Value = "  +3.98%  "

# row 9 - Solana
Set-TextValue "D9" "24.11"
$ws.Range("E9").Value = "  -0.24%  "

# row 10 - Cardano
$ws.Range("E10").Value = "  -1.62%  "

# row 11 - Dogecoin
$ws.Range("E11").Value = "  -1.57%  "

# row 12 - TRON
Set-TextValue "D12" "0.0889"
$ws.Range("E12").Value = "  -0.06%  "

# row 14 - WrappedEther
Set-TextValue "D14" "1.576.77"
$ws.Range("E14").Value = "  -0.82%  "

# row 15 - Polygon
$ws.Range("E15").Value = "  -1.86%  "

# row 16 - WrappedBTC
Set-TextValue "D16" "28.532.71"
$ws.Range("E16").Value = "  +0.78%  "

# row 18 - Litecoin
Set-TextValue "D18" "62.32"
$ws.Range("E18").Value = "  -1.36%  "

# row 19 - BitcoinCash
Set-TextValue "D19" "230.45"
$ws.Range("E19").Value = "  +1.09%  "

# row 20 - Chainlink
$ws.Range("E20").Value = "  -1.55%  "

# row 21 - ShibaInu
$ws.Range("E21").Value = "  -2.69%  "

# row 23 - Uniswap
Set-TextValue "D23" "3.89"
$ws.Range("E23").Value = "  -5.30%  "

# row 24 - Avalanche
Set-TextValue "D24" "9.12"
$ws.Range("E24").Value = "  -2.37%  "

# row 25 - Toncoin
$ws.Range("E25").Value = "  +9.99%  "

# row 26 - Monero
Set-TextValue "D26" "151.82"
$ws.Range("E26").Value = "  -0.04%  "

# row 27 - EthereumClassic
Set-TextValue "D27" "15.01"
$ws.Range("E27").Value = "  -1.39%  "

# row 28 - Cosmos
$ws.Range("E28").Value = "  -2.32%  "

# row 29 - Stellar
$ws.Range("E29").Value = "  -3.22%  "

# row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.05%  "

# row 31 - Hedera
Set-TextValue "D31" "0.0485"
$ws.Range("E31").Value = "  +2.28%  "

# row 32 - PancakeSwap
$ws.Range("E32").Value = "  -2.68%  "

# row 33 - Filecoin
Set-TextValue "D33" "3.21"
$ws.Range("E33").Value = "  -1.01%  "

# row 34 - InternetComputer(DFINITY)
Set-TextValue "D34" "3.10"
$ws.Range("E34").Value = "  -1.77%  "

# row 35 - Maker
Set-TextValue "D35" "1.390.63"
$ws.Range("E35").Value = "  -0.94%  "

# row 36 - TrustWalletToken
$ws.Range("E36").Value = "  +1.49%  "

# row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -3.41%  "

# row 38 - HuobiToken
$ws.Range("E38").Value = "  +0.95%  "

# row 39 - MXToken
$ws.Range("E39").Value = "  +3.10%  "

# row 40 - VeChain
$ws.Range("E40").Value = "  -0.70%  "

# row 41 - ImmutableX
$ws.Range("E41").Value = "  -3.04%  "

# row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.08%  "

# row 43 - RenderToken
$ws.Range("E43").Value = "  +0.61%  "

# row 44 - ARBITRUM
$ws.Range("E44").Value = "  -3.00%  "

# row 45 - Kaspa
$ws.Range("E45").Value = "  +2.25%  "

# row 46 - FraxShare
$ws.Range("E46").Value = "  -2.71%  "

# row 47 - WEMIXToken
Set-TextValue "D47" "0.967"
$ws.Range("E47").Value = "  -1.92%  "

# row 48 - Aave
Set-TextValue "D48" "63.00"
$ws.Range("E48").Value = "  -2.11%  "

# row 49 - RocketPoolETH
Set-TextValue "D49" "1.708.49"
$ws.Range("E49").Value = "  -1.46%  "

# row 50 - Quant
Set-TextValue "D50" "86.49"
$ws.Range("E50").Value = "  -1.37%  "

# row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -0.84%  "
